# Update the jxls demo report template placeholders to use lowercase
# "row.*" property names instead of the old uppercase ones, and update
# the active selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A10").Value = '${row.order_id}'
$ws.Range("B10").Value = '${row.city_name}'
$ws.Range("C10").Value = '${row.item_name}'
$ws.Range("D10").Value = '${row.order_date}'
$ws.Range("E10").Value = '${row.volume}'

$ws.Activate()
$ws.Range("F11").Select()
